# Sprint 2 header-table update for mailer.html.erb code review doc:
#   Sprint No.  : 1          -> 2
#   Review Date : 02/09/18   -> 02/21/18
#
# The "Sprint No." / "Review Date" values live in the first table of the
# document (the review metadata table). Address the exact cells directly
# rather than a document-wide Find/Replace so that unrelated occurrences
# of "1" elsewhere in the checklist (e.g. "1.  License", dates, etc.)
# are left untouched.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# A table cell's Range.Text carries a trailing cell-mark (CR + BEL); strip
# those control characters off before comparing the visible content.
$cellMarkPattern = "[\x07\x0d]+$"

# Row 2, Column 4 holds the Sprint No. value ("1").
$sprintCell = $table.Cell(2, 4)
$sprintText = $sprintCell.Range.Text -replace $cellMarkPattern, ""
if ($sprintText -eq "1") {
    $sprintCell.Range.Text = "2"
}

# Row 3, Column 2 holds the Review Date value ("02/09/18").
$dateCell = $table.Cell(3, 2)
$dateText = $dateCell.Range.Text -replace $cellMarkPattern, ""
if ($dateText -eq "02/09/18") {
    $dateCell.Range.Text = "02/21/18"
}
